$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.689.43"
$ws.Range("E2").Value = "  -0.49%  "
$ws.Range("D3").Value = "'2.298.16"
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'301.17"
$ws.Range("E5").Value = "  -1.48%  "
$ws.Range("D6").Value = "'95.99"
$ws.Range("E6").Value = "  -1.17%  "
$ws.Range("D7").Value = "'0.512"
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "'0.495"
$ws.Range("E9").Value = "  -1.68%  "
$ws.Range("D10").Value = "'34.77"
$ws.Range("E10").Value = "  -2.36%  "
$ws.Range("E11").Value = "  +5.07%  "
$ws.Range("D12").Value = "'0.0788"
$ws.Range("E12").Value = "  -0.90%  "
$ws.Range("E13").Value = "  -0.26%  "
$ws.Range("E14").Value = "  +0.42%  "
$ws.Range("D15").Value = "'2.650.20"
$ws.Range("E15").Value = "  -0.36%  "
$ws.Range("D16").Value = "'2.298.78"
$ws.Range("E16").Value = "  +0.02%  "
$ws.Range("E17").Value = "  +0.25%  "
$ws.Range("D18").Value = "'42.621.51"
$ws.Range("D19").Value = "'12.32"
$ws.Range("E19").Value = "  -6.62%  "
$ws.Range("E20").Value = "  -1.46%  "
$ws.Range("E21").Value = "  -0.28%  "
$ws.Range("D22").Value = "'67.86"
$ws.Range("E22").Value = "  +0.44%  "
$ws.Range("E23").Value = "  +3.22%  "
$ws.Range("D24").Value = "'235.37"
$ws.Range("E24").Value = "  -0.54%  "
$ws.Range("E25").Value = "  +0.17%  "
$ws.Range("D26").Value = "'2.41"
$ws.Range("E26").Value = "  -2.36%  "
$ws.Range("D27").Value = "'24.62"
$ws.Range("E27").Value = "  -3.46%  "
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").Value = "'2.06"
$ws.Range("E28").Value = "  -0.32%  "
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "'164.84"
$ws.Range("E29").Value = "  -1.48%  "
$ws.Range("D30").Value = "'9.07"
$ws.Range("E30").Value = "  -0.23%  "
$ws.Range("D31").Value = "'32.18"
$ws.Range("E31").Value = "  -2.53%  "
$ws.Range("E32").Value = "  -0.04%  "
$ws.Range("D33").Value = "'4.97"
$ws.Range("E33").Value = "  -0.73%  "
$ws.Range("D34").Value = "'17.55"
$ws.Range("E34").Value = "  +1.00%  "
$ws.Range("D35").Value = "'4.45"
$ws.Range("E35").Value = "  -7.27%  "
$ws.Range("D36").Value = "'0.0702"
$ws.Range("E36").Value = "  +1.48%  "
$ws.Range("E37").Value = "  -2.99%  "
$ws.Range("E39").Value = "  -0.22%  "
$ws.Range("E40").Value = "  -0.19%  "
$ws.Range("D41").Value = "'0.109"
$ws.Range("E41").Value = "  -1.04%  "
$ws.Range("D42").Value = "'19.81"
$ws.Range("E42").Value = "  +7.58%  "
$ws.Range("D43").Value = "'1.970.49"
$ws.Range("E43").Value = "  -1.95%  "
$ws.Range("E44").Value = "  +4.49%  "
$ws.Range("E45").Value = "  -1.03%  "
$ws.Range("E46").Value = "  -3.19%  "
$ws.Range("D47").Value = "'2.76"
$ws.Range("E47").Value = "  -0.85%  "
$ws.Range("D48").Value = "'2.90"
$ws.Range("E48").Value = "  -1.18%  "
$ws.Range("D49").Value = "'2.523.42"
$ws.Range("E49").Value = "  -0.17%  "
$ws.Range("D50").Value = "'53.10"
$ws.Range("E50").Value = "  -1.60%  "
$ws.Range("D51").Value = "'71.59"
$ws.Range("E51").Value = "  -0.55%  "
